# "test RAG with many model"
# - Clears the old single-row summary line from "Case 1" and rebuilds the
#   table with extra MODEL EMBEDDING / MODEL GEN ANSWER / Status columns
#   (plus a stray note in J3).
# - Adds a second worksheet ("Sheet1") with a small model/max-length/
#   hidden-state-shape comparison table, and makes it the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# "Case 1" sheet: wipe the old row 1 ("MODEL = gpt-3.5-turbo" in F1) and
# the old row 2/3/4 layout, then re-enter the new table.
# ---------------------------------------------------------------------
$ws.Range("A1:F1").ClearContents()

# Header row (row 2)
$ws.Range("A2").Value = "stt"
$ws.Range("B2").Value = "CHUNK SIZE"
$ws.Range("C2").Value = "CHUNK OVERLAP"
$ws.Range("D2").Value = "MODEL EMBEDDING"
$ws.Range("E2").Value = "MODEL GEN ANSWER"
$ws.Range("F2").Value = "Status"

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 256
$ws.Range("C3").Value = 20
$ws.Range("D3").Value = "`ntext-embedding-ada-002 (from openai)"
$ws.Range("D3").WrapText = $true
$ws.Range("E3").Value = " gpt-3.5-turbo"
$ws.Range("F3").Value = "Done"
$ws.Range("J3").Value = "phobert(base): unefficient"
$ws.Range("J3").WrapText = $true

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 2000
$ws.Range("C4").Value = 200
$ws.Range("D4").Value = "`ntext-embedding-ada-002 (from openai)"
$ws.Range("D4").WrapText = $true

# Row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 1024
$ws.Range("C5").Value = 100
$ws.Range("D5").Value = "`ntext-embedding-ada-002 (from openai)"
$ws.Range("D5").WrapText = $true
$ws.Range("E5").Value = "gpt-3.5-turbo"
$ws.Range("F5").Value = "Done"

# Row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 256
$ws.Range("C6").Value = 20
$ws.Range("D6").Value = "phobert(base)"
$ws.Range("D6").WrapText = $true
$ws.Range("E6").Value = "gpt-3.5-turbo"

# Row 7
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = 256
$ws.Range("C7").Value = 20
$ws.Range("D7").Value = "phobert(base)"
$ws.Range("D7").WrapText = $true
$ws.Range("E7").Value = "mistralai/Mixtral-8x7B-Instruct-v0.1"
$ws.Range("E7").Font.Name = "Consolas"
$ws.Range("E7").VerticalAlignment = -4108

# Row 8
$ws.Range("D8").Value = "vietnamese sbert"
$ws.Range("D8").Font.Name = "Segoe UI Historic"
$ws.Range("D8").Font.Color = 328965

# Row heights to match the wrapped, multi-line content
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 45
$ws.Rows.Item(5).RowHeight = 45
$ws.Rows.Item(8).RowHeight = 16.5

# Column widths for the columns that grew to fit the longer text
$ws.Columns.Item(4).ColumnWidth = 39.5
$ws.Columns.Item(5).ColumnWidth = 41.833333333333336
$ws.Columns.Item(10).ColumnWidth = 29.5

$ws.PageSetup.Orientation = 1

$ws.Range("E21").Select() | Out-Null

# ---------------------------------------------------------------------
# New "Sheet1" worksheet, placed right after "Case 1" and left as the
# active tab (matches the new activeTab="1" on the workbook).
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws)
$ws2.Name = "Sheet1"

$ws2.Range("A1").Value = "Model name"
$ws2.Range("B1").Value = "Max length"
$ws2.Range("C1").Value = "Last_hidden_state shape"

$ws2.Range("A2").Value = "phobert"
$ws2.Range("B2").Value = 256
$ws2.Range("C2").Value = 768

$ws2.Range("A3").Value = "bartpho"
$ws2.Range("B3").Value = 1024
$ws2.Range("C3").Value = 1024

$ws2.Columns.Item(1).ColumnWidth = 25.666666666666668
$ws2.Columns.Item(2).ColumnWidth = 29.666666666666668
$ws2.Columns.Item(3).ColumnWidth = 27.0

$ws2.Range("F14").Select() | Out-Null
